$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 126
$ws.Cells.Item(2, 3).Value = "flower/flower026.jpg"
$ws.Cells.Item(2, 4).Value = "pflegen"
$ws.Cells.Item(2, 5).Value = "flower"
$ws.Cells.Item(3, 2).Value = 66
$ws.Cells.Item(3, 3).Value = "dog/dog001.jpg"
$ws.Cells.Item(3, 4).Value = "mieten"
$ws.Cells.Item(3, 5).Value = "dog"
$ws.Cells.Item(4, 2).Value = 85
$ws.Cells.Item(4, 3).Value = "flower/flower016.jpg"
$ws.Cells.Item(4, 4).Value = "schenken"
$ws.Cells.Item(4, 5).Value = "flower"
$ws.Cells.Item(5, 2).Value = 50
$ws.Cells.Item(5, 3).Value = "flower/flower002.jpg"
$ws.Cells.Item(5, 4).Value = "stärken"
$ws.Cells.Item(5, 5).Value = "flower"
$ws.Cells.Item(6, 2).Value = 41
$ws.Cells.Item(6, 3).Value = "dog/dog013.jpg"
$ws.Cells.Item(6, 4).Value = "dauern"
$ws.Cells.Item(6, 5).Value = "dog"
$ws.Cells.Item(7, 2).Value = 84
$ws.Cells.Item(7, 3).Value = "dog/dog008.jpg"
$ws.Cells.Item(7, 4).Value = "haken"
$ws.Cells.Item(7, 5).Value = "dog"
$ws.Cells.Item(8, 2).Value = 77
$ws.Cells.Item(8, 3).Value = "flower/flower031.jpg"
$ws.Cells.Item(8, 4).Value = "wiegen"
$ws.Cells.Item(8, 5).Value = "flower"
$ws.Cells.Item(9, 2).Value = 121
$ws.Cells.Item(9, 3).Value = "flower/flower027.jpg"
$ws.Cells.Item(9, 4).Value = "lehnen"
$ws.Cells.Item(9, 5).Value = "flower"
$ws.Cells.Item(10, 2).Value = 115
$ws.Cells.Item(10, 3).Value = "dog/dog029.jpg"
$ws.Cells.Item(10, 4).Value = "fesseln"
$ws.Cells.Item(10, 5).Value = "dog"
$ws.Cells.Item(11, 2).Value = 114
$ws.Cells.Item(11, 3).Value = "dog/dog002.jpg"
$ws.Cells.Item(11, 4).Value = "gründen"
$ws.Cells.Item(11, 5).Value = "dog"
$ws.Cells.Item(12, 2).Value = 57
$ws.Cells.Item(12, 3).Value = "dog/dog019.jpg"
$ws.Cells.Item(12, 4).Value = "währen"
$ws.Cells.Item(12, 5).Value = "dog"
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(13, 3).Value = "flower/flower003.jpg"
$ws.Cells.Item(13, 4).Value = "jubeln"
$ws.Cells.Item(13, 5).Value = "flower"
$ws.Cells.Item(14, 2).Value = 55
$ws.Cells.Item(14, 3).Value = "dog/dog011.jpg"
$ws.Cells.Item(14, 4).Value = "drohen"
$ws.Cells.Item(14, 5).Value = "dog"
$ws.Cells.Item(15, 2).Value = 104
$ws.Cells.Item(15, 3).Value = "flower/flower006.jpg"
$ws.Cells.Item(15, 4).Value = "laufen"
$ws.Cells.Item(15, 5).Value = "flower"
$ws.Cells.Item(16, 2).Value = 52
$ws.Cells.Item(16, 3).Value = "dog/dog005.jpg"
$ws.Cells.Item(16, 4).Value = "fühlen"
$ws.Cells.Item(16, 5).Value = "dog"
$ws.Cells.Item(17, 2).Value = 37
$ws.Cells.Item(17, 3).Value = "flower/flower007.jpg"
$ws.Cells.Item(17, 4).Value = "bleiben"
$ws.Cells.Item(17, 5).Value = "flower"
$ws.Cells.Item(18, 2).Value = 96
$ws.Cells.Item(18, 3).Value = "flower/flower005.jpg"
$ws.Cells.Item(18, 4).Value = "strahlen"
$ws.Cells.Item(18, 5).Value = "flower"
$ws.Cells.Item(19, 2).Value = 42
$ws.Cells.Item(19, 3).Value = "dog/dog027.jpg"
$ws.Cells.Item(19, 4).Value = "rasen"
$ws.Cells.Item(19, 5).Value = "dog"
$ws.Cells.Item(20, 2).Value = 19
$ws.Cells.Item(20, 3).Value = "dog/dog014.jpg"
$ws.Cells.Item(20, 4).Value = "runden"
$ws.Cells.Item(20, 5).Value = "dog"
$ws.Cells.Item(21, 2).Value = 103
$ws.Cells.Item(21, 3).Value = "flower/flower018.jpg"
$ws.Cells.Item(21, 4).Value = "kehren"
$ws.Cells.Item(21, 5).Value = "flower"
$ws.Cells.Item(22, 2).Value = 65
$ws.Cells.Item(22, 3).Value = "dog/dog003.jpg"
$ws.Cells.Item(22, 4).Value = "hauen"
$ws.Cells.Item(22, 5).Value = "dog"
$ws.Cells.Item(23, 2).Value = 109
$ws.Cells.Item(23, 3).Value = "dog/dog021.jpg"
$ws.Cells.Item(23, 4).Value = "hoffen"
$ws.Cells.Item(23, 5).Value = "dog"
$ws.Cells.Item(24, 2).Value = 27
$ws.Cells.Item(24, 3).Value = "dog/dog017.jpg"
$ws.Cells.Item(24, 4).Value = "fliehen"
$ws.Cells.Item(24, 5).Value = "dog"
$ws.Cells.Item(25, 2).Value = 98
$ws.Cells.Item(25, 3).Value = "dog/dog012.jpg"
$ws.Cells.Item(25, 4).Value = "sondern"
$ws.Cells.Item(25, 5).Value = "dog"
$ws.Cells.Item(26, 2).Value = 75
$ws.Cells.Item(26, 3).Value = "flower/flower013.jpg"
$ws.Cells.Item(26, 4).Value = "backen"
$ws.Cells.Item(26, 5).Value = "flower"
$ws.Cells.Item(27, 2).Value = 15
$ws.Cells.Item(27, 3).Value = "dog/dog020.jpg"
$ws.Cells.Item(27, 4).Value = "langen"
$ws.Cells.Item(27, 5).Value = "dog"
$ws.Cells.Item(28, 2).Value = 5
$ws.Cells.Item(28, 3).Value = "flower/flower025.jpg"
$ws.Cells.Item(28, 4).Value = "klappen"
$ws.Cells.Item(28, 5).Value = "flower"
$ws.Cells.Item(29, 2).Value = 40
$ws.Cells.Item(29, 3).Value = "flower/flower015.jpg"
$ws.Cells.Item(29, 4).Value = "gelten"
$ws.Cells.Item(29, 5).Value = "flower"
$ws.Cells.Item(30, 2).Value = 54
$ws.Cells.Item(30, 3).Value = "flower/flower010.jpg"
$ws.Cells.Item(30, 4).Value = "tauschen"
$ws.Cells.Item(30, 5).Value = "flower"
$ws.Cells.Item(31, 2).Value = 24
$ws.Cells.Item(31, 3).Value = "flower/flower001.jpg"
$ws.Cells.Item(31, 4).Value = "husten"
$ws.Cells.Item(31, 5).Value = "flower"
$ws.Cells.Item(32, 2).Value = 70
$ws.Cells.Item(32, 3).Value = "dog/dog004.jpg"
$ws.Cells.Item(32, 4).Value = "saufen"
$ws.Cells.Item(32, 5).Value = "dog"
$ws.Cells.Item(33, 2).Value = 125
$ws.Cells.Item(33, 3).Value = "flower/flower014.jpg"
$ws.Cells.Item(33, 4).Value = "antun"
$ws.Cells.Item(33, 5).Value = "flower"
